# task-list.xlsx -- "item update with tax"
# 1) Sets explicit row heights on a batch of existing rows (taller rows for
#    wrapped multi-line task descriptions).
# 2) Appends 6 new task rows (314-319) with their own row heights.
# 3) Resets the sheet view to the top-left / A1 selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Row-height adjustments on existing rows
# ---------------------------------------------------------------------------
$rowHeights = @(
    @(278, 60),
    @(279, 45),
    @(280, 30),
    @(282, 45),
    @(283, 30),
    @(285, 30),
    @(286, 30),
    @(287, 30),
    @(288, 45),
    @(291, 45),
    @(300, 30),
    @(301, 30),
    @(302, 30),
    @(304, 30),
    @(305, 30),
    @(306, 90),
    @(308, 45),
    @(313, 30)
)

foreach ($pair in $rowHeights) {
    $rowNum = $pair[0]
    $height = $pair[1]
    $ws.Rows.Item($rowNum).RowHeight = $height
}

# ---------------------------------------------------------------------------
# 2) Append six new rows of task data (rows 314-319)
# ---------------------------------------------------------------------------
$newTasks = @(
    "Sup. ret. weighted avg confirmation/",
    "partial grn from same po/",
    "ap invoice batch testing/",
    "sequence+with outletid in sales\pay/-db sequence used",
    "po screen - re-named column for ship to location instead of dff/",
    "identify warehouse -> outlet transfer from a notification sent based on the inter-stock and Hold payment till inter-stock doneX"
)

$newRowHeights = @{
    317 = 30
    318 = 30
    319 = 45
}

$startRow = 314
for ($i = 0; $i -lt $newTasks.Length; $i++) {
    $rowNum = $startRow + $i

    $ws.Cells.Item($rowNum, 1).Value = ($rowNum - 1)

    $bCell = $ws.Cells.Item($rowNum, 2)
    $bCell.Value = $newTasks[$i]
    $bCell.WrapText = $true

    $cCell = $ws.Cells.Item($rowNum, 3)
    $cCell.Value = 810626
    $cCell.NumberFormat = "mmm-yy"

    $dCell = $ws.Cells.Item($rowNum, 4)
    $dCell.Value = 810626
    $dCell.NumberFormat = "mmm-yy"

    if ($newRowHeights.ContainsKey($rowNum)) {
        $ws.Rows.Item($rowNum).RowHeight = $newRowHeights[$rowNum]
    }
}

# ---------------------------------------------------------------------------
# 3) Reset the view: scroll back to the top and select A1
# ---------------------------------------------------------------------------
$ws.Range("A1").Select()
